$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LCGA")

# ---------------------------------------------------------------------------
# Header row (row 1) -- rename existing "ok" columns to be flexmix-specific,
# and add three new "custom flexmix" columns plus a free-text "note" column.
# ---------------------------------------------------------------------------
$ws.Range("E1").Value = "betas ok (flexmix)"
$ws.Range("F1").Value = "Rs ok (flexmix)"
$ws.Range("G1").Value = "pis ok (flexmix)"
$ws.Range("H1").Value = "betas ok (custom flexmix)"
$ws.Range("I1").Value = "Rs ok custom flexmix)"
$ws.Range("J1").Value = "pis ok (custom flexmix)"
$ws.Range("K1").Value = "note"

# ---------------------------------------------------------------------------
# Row 2
# ---------------------------------------------------------------------------
$ws.Range("H2").Value = "ok"
$ws.Range("I2").Value = "ok"
$ws.Range("J2").Value = "ok"

# ---------------------------------------------------------------------------
# Row 3
# ---------------------------------------------------------------------------
$ws.Range("H3").Value = "ok"
$ws.Range("I3").Value = "ok"
$ws.Range("J3").Value = "ok"
$ws.Range("K3").Value = "needed to rerun flexmix to have same res!"

# ---------------------------------------------------------------------------
# Row 4
# ---------------------------------------------------------------------------
$ws.Range("H4").Value = "almost (3 cases)"
$ws.Range("I4").Value = "almost (3 cases)"
$ws.Range("J4").Value = "almost (3 cases)"

# ---------------------------------------------------------------------------
# Row 5 -- also gets a note with mixed (partially bold) rich text.
# ---------------------------------------------------------------------------
$ws.Range("H5").Value = "no"
$ws.Range("I5").Value = "no"
$ws.Range("J5").Value = "no"

$note = "in our simulation, some iterations found a result very close to flexmix, but duer to mulststart the final output is different, which corroborates that our algorithm was better than flexmix in this dataset"
$ws.Range("K5").Value = $note
$boldStart = $note.IndexOf("better") + 1
$boldLen = 6
$ws.Range("K5").Characters($boldStart, $boldLen).Font.Bold = $true
$afterStart = $boldStart + $boldLen
$afterLen = $note.Length - $afterStart + 1
$ws.Range("K5").Characters($afterStart, $afterLen).Font.Bold = $false
$ws.Range("K5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 63.6

# ---------------------------------------------------------------------------
# Row 6
# ---------------------------------------------------------------------------
$ws.Range("H6").Value = "almost (3 cases)"
$ws.Range("I6").Value = "almost (3 cases)"
$ws.Range("J6").Value = "almost (4 cases)"

# ---------------------------------------------------------------------------
# Column sizing to roughly match the widened layout (columns E..K).
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 15.1666666666667
$ws.Columns.Item(6).ColumnWidth = 12.5
$ws.Columns.Item(7).ColumnWidth = 13
$ws.Columns.Item(8).ColumnWidth = 22
$ws.Columns.Item(9).ColumnWidth = 20.1666666666667
$ws.Columns.Item(10).ColumnWidth = 19.6666666666667
$ws.Columns.Item(11).ColumnWidth = 42.3333333333333

# ---------------------------------------------------------------------------
# Selection moves one row below the new last data row, at the new last column.
# ---------------------------------------------------------------------------
[void]$ws.Range("K7").Select()
